$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 284-303 and append new rows 304-305
# Data reconstructed from the target diff (rows shift down by 2; two new rows
# of data are inserted at 284-285; two new rows are appended at 304-305).

# Row 284
$ws.Cells.Item(284,1).Value = 10
$ws.Cells.Item(284,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(284,3).Value = 'La Araucanía'
$ws.Cells.Item(284,4).Value = 45106
$ws.Cells.Item(284,5).Value = 9
$ws.Cells.Item(284,6).Value = 100112013
$ws.Cells.Item(284,7).Value = 'Alcachofa'
$ws.Cells.Item(284,8).Value = 'Española'
$ws.Cells.Item(284,9).Value = 'Primera'
$ws.Cells.Item(284,10).Value = 450
$ws.Cells.Item(284,11).Value = 500
$ws.Cells.Item(284,12).Value = 550
$ws.Cells.Item(284,13).Value = 528
$ws.Cells.Item(284,14).Value = '$/unidad'
$ws.Cells.Item(284,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(284,16).Value = 528
$ws.Cells.Item(284,17).Value = 1
$ws.Cells.Item(284,18).Value = 'Hortaliza'

# Row 285
$ws.Cells.Item(285,1).Value = 10
$ws.Cells.Item(285,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(285,3).Value = 'La Araucanía'
$ws.Cells.Item(285,4).Value = 45106
$ws.Cells.Item(285,5).Value = 9
$ws.Cells.Item(285,6).Value = 100112013
$ws.Cells.Item(285,7).Value = 'Alcachofa'
$ws.Cells.Item(285,8).Value = 'Madrigal'
$ws.Cells.Item(285,9).Value = 'Extra'
$ws.Cells.Item(285,10).Value = 155
$ws.Cells.Item(285,11).Value = 16000
$ws.Cells.Item(285,12).Value = 16000
$ws.Cells.Item(285,13).Value = 16000
$ws.Cells.Item(285,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(285,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(285,16).Value = 533
$ws.Cells.Item(285,17).Value = 30
$ws.Cells.Item(285,18).Value = 'Hortaliza'

# Row 286
$ws.Cells.Item(286,1).Value = 10
$ws.Cells.Item(286,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(286,3).Value = 'La Araucanía'
$ws.Cells.Item(286,4).Value = 44383
$ws.Cells.Item(286,5).Value = 9
$ws.Cells.Item(286,6).Value = 100112013
$ws.Cells.Item(286,7).Value = 'Alcachofa'
$ws.Cells.Item(286,8).Value = 'Española'
$ws.Cells.Item(286,9).Value = 'Primera'
$ws.Cells.Item(286,10).Value = 50
$ws.Cells.Item(286,11).Value = 18000
$ws.Cells.Item(286,12).Value = 18000
$ws.Cells.Item(286,13).Value = 18000
$ws.Cells.Item(286,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(286,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(286,16).Value = 600
$ws.Cells.Item(286,17).Value = 30
$ws.Cells.Item(286,18).Value = 'Hortaliza'

# Row 287
$ws.Cells.Item(287,1).Value = 10
$ws.Cells.Item(287,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(287,3).Value = 'La Araucanía'
$ws.Cells.Item(287,4).Value = 44383
$ws.Cells.Item(287,5).Value = 9
$ws.Cells.Item(287,6).Value = 100112013
$ws.Cells.Item(287,7).Value = 'Alcachofa'
$ws.Cells.Item(287,8).Value = 'Madrigal'
$ws.Cells.Item(287,9).Value = 'Primera'
$ws.Cells.Item(287,10).Value = 70
$ws.Cells.Item(287,11).Value = 18000
$ws.Cells.Item(287,12).Value = 18000
$ws.Cells.Item(287,13).Value = 18000
$ws.Cells.Item(287,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(287,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(287,16).Value = 450
$ws.Cells.Item(287,17).Value = 40
$ws.Cells.Item(287,18).Value = 'Hortaliza'

# Row 288
$ws.Cells.Item(288,1).Value = 10
$ws.Cells.Item(288,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(288,3).Value = 'La Araucanía'
$ws.Cells.Item(288,4).Value = 45068
$ws.Cells.Item(288,5).Value = 9
$ws.Cells.Item(288,6).Value = 100112013
$ws.Cells.Item(288,7).Value = 'Alcachofa'
$ws.Cells.Item(288,8).Value = 'Madrigal'
$ws.Cells.Item(288,9).Value = 'Primera'
$ws.Cells.Item(288,10).Value = 140
$ws.Cells.Item(288,11).Value = 20000
$ws.Cells.Item(288,12).Value = 20000
$ws.Cells.Item(288,13).Value = 20000
$ws.Cells.Item(288,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(288,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(288,16).Value = 500
$ws.Cells.Item(288,17).Value = 40
$ws.Cells.Item(288,18).Value = 'Hortaliza'

# Row 289
$ws.Cells.Item(289,1).Value = 10
$ws.Cells.Item(289,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(289,3).Value = 'La Araucanía'
$ws.Cells.Item(289,4).Value = 44756
$ws.Cells.Item(289,5).Value = 9
$ws.Cells.Item(289,6).Value = 100112013
$ws.Cells.Item(289,7).Value = 'Alcachofa'
$ws.Cells.Item(289,8).Value = 'Madrigal'
$ws.Cells.Item(289,9).Value = 'Extra'
$ws.Cells.Item(289,10).Value = 80
$ws.Cells.Item(289,11).Value = 22000
$ws.Cells.Item(289,12).Value = 22000
$ws.Cells.Item(289,13).Value = 22000
$ws.Cells.Item(289,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(289,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(289,16).Value = 733
$ws.Cells.Item(289,17).Value = 30
$ws.Cells.Item(289,18).Value = 'Hortaliza'

# Row 290
$ws.Cells.Item(290,1).Value = 10
$ws.Cells.Item(290,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(290,3).Value = 'La Araucanía'
$ws.Cells.Item(290,4).Value = 44756
$ws.Cells.Item(290,5).Value = 9
$ws.Cells.Item(290,6).Value = 100112013
$ws.Cells.Item(290,7).Value = 'Alcachofa'
$ws.Cells.Item(290,8).Value = 'Madrigal'
$ws.Cells.Item(290,9).Value = 'Primera'
$ws.Cells.Item(290,10).Value = 450
$ws.Cells.Item(290,11).Value = 18000
$ws.Cells.Item(290,12).Value = 19000
$ws.Cells.Item(290,13).Value = 18556
$ws.Cells.Item(290,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(290,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(290,16).Value = 464
$ws.Cells.Item(290,17).Value = 40
$ws.Cells.Item(290,18).Value = 'Hortaliza'

# Row 291
$ws.Cells.Item(291,1).Value = 10
$ws.Cells.Item(291,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(291,3).Value = 'La Araucanía'
$ws.Cells.Item(291,4).Value = 44810
$ws.Cells.Item(291,5).Value = 9
$ws.Cells.Item(291,6).Value = 100112013
$ws.Cells.Item(291,7).Value = 'Alcachofa'
$ws.Cells.Item(291,8).Value = 'Madrigal'
$ws.Cells.Item(291,9).Value = 'Primera'
$ws.Cells.Item(291,10).Value = 150
$ws.Cells.Item(291,11).Value = 12000
$ws.Cells.Item(291,12).Value = 12000
$ws.Cells.Item(291,13).Value = 12000
$ws.Cells.Item(291,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(291,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(291,16).Value = 300
$ws.Cells.Item(291,17).Value = 40
$ws.Cells.Item(291,18).Value = 'Hortaliza'

# Row 292
$ws.Cells.Item(292,1).Value = 10
$ws.Cells.Item(292,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(292,3).Value = 'La Araucanía'
$ws.Cells.Item(292,4).Value = 44442
$ws.Cells.Item(292,5).Value = 9
$ws.Cells.Item(292,6).Value = 100112013
$ws.Cells.Item(292,7).Value = 'Alcachofa'
$ws.Cells.Item(292,8).Value = 'Argentina(o)'
$ws.Cells.Item(292,9).Value = 'Primera'
$ws.Cells.Item(292,10).Value = 90
$ws.Cells.Item(292,11).Value = 13000
$ws.Cells.Item(292,12).Value = 13000
$ws.Cells.Item(292,13).Value = 13000
$ws.Cells.Item(292,14).Value = '$/caja 50 unidades'
$ws.Cells.Item(292,15).Value = 'Región Metropolitana'
$ws.Cells.Item(292,16).Value = 260
$ws.Cells.Item(292,17).Value = 50
$ws.Cells.Item(292,18).Value = 'Hortaliza'

# Row 293
$ws.Cells.Item(293,1).Value = 10
$ws.Cells.Item(293,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(293,3).Value = 'La Araucanía'
$ws.Cells.Item(293,4).Value = 44504
$ws.Cells.Item(293,5).Value = 9
$ws.Cells.Item(293,6).Value = 100112013
$ws.Cells.Item(293,7).Value = 'Alcachofa'
$ws.Cells.Item(293,8).Value = 'Madrigal'
$ws.Cells.Item(293,9).Value = 'Primera'
$ws.Cells.Item(293,10).Value = 95
$ws.Cells.Item(293,11).Value = 12000
$ws.Cells.Item(293,12).Value = 12000
$ws.Cells.Item(293,13).Value = 12000
$ws.Cells.Item(293,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(293,15).Value = 'Región del Maule'
$ws.Cells.Item(293,16).Value = 300
$ws.Cells.Item(293,17).Value = 40
$ws.Cells.Item(293,18).Value = 'Hortaliza'

# Row 294
$ws.Cells.Item(294,1).Value = 10
$ws.Cells.Item(294,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(294,3).Value = 'La Araucanía'
$ws.Cells.Item(294,4).Value = 44769
$ws.Cells.Item(294,5).Value = 9
$ws.Cells.Item(294,6).Value = 100112013
$ws.Cells.Item(294,7).Value = 'Alcachofa'
$ws.Cells.Item(294,8).Value = 'Española'
$ws.Cells.Item(294,9).Value = 'Primera'
$ws.Cells.Item(294,10).Value = 120
$ws.Cells.Item(294,11).Value = 18000
$ws.Cells.Item(294,12).Value = 18000
$ws.Cells.Item(294,13).Value = 18000
$ws.Cells.Item(294,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(294,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(294,16).Value = 600
$ws.Cells.Item(294,17).Value = 30
$ws.Cells.Item(294,18).Value = 'Hortaliza'

# Row 295
$ws.Cells.Item(295,1).Value = 10
$ws.Cells.Item(295,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(295,3).Value = 'La Araucanía'
$ws.Cells.Item(295,4).Value = 44769
$ws.Cells.Item(295,5).Value = 9
$ws.Cells.Item(295,6).Value = 100112013
$ws.Cells.Item(295,7).Value = 'Alcachofa'
$ws.Cells.Item(295,8).Value = 'Madrigal'
$ws.Cells.Item(295,9).Value = 'Primera'
$ws.Cells.Item(295,10).Value = 200
$ws.Cells.Item(295,11).Value = 15000
$ws.Cells.Item(295,12).Value = 15000
$ws.Cells.Item(295,13).Value = 15000
$ws.Cells.Item(295,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(295,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(295,16).Value = 375
$ws.Cells.Item(295,17).Value = 40
$ws.Cells.Item(295,18).Value = 'Hortaliza'

# Row 296
$ws.Cells.Item(296,1).Value = 10
$ws.Cells.Item(296,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(296,3).Value = 'La Araucanía'
$ws.Cells.Item(296,4).Value = 45104
$ws.Cells.Item(296,5).Value = 9
$ws.Cells.Item(296,6).Value = 100112013
$ws.Cells.Item(296,7).Value = 'Alcachofa'
$ws.Cells.Item(296,8).Value = 'Madrigal'
$ws.Cells.Item(296,9).Value = 'Extra'
$ws.Cells.Item(296,10).Value = 110
$ws.Cells.Item(296,11).Value = 16000
$ws.Cells.Item(296,12).Value = 16000
$ws.Cells.Item(296,13).Value = 16000
$ws.Cells.Item(296,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(296,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(296,16).Value = 533
$ws.Cells.Item(296,17).Value = 30
$ws.Cells.Item(296,18).Value = 'Hortaliza'

# Row 297
$ws.Cells.Item(297,1).Value = 10
$ws.Cells.Item(297,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(297,3).Value = 'La Araucanía'
$ws.Cells.Item(297,4).Value = 45096
$ws.Cells.Item(297,5).Value = 9
$ws.Cells.Item(297,6).Value = 100112013
$ws.Cells.Item(297,7).Value = 'Alcachofa'
$ws.Cells.Item(297,8).Value = 'Madrigal'
$ws.Cells.Item(297,9).Value = 'Primera'
$ws.Cells.Item(297,10).Value = 800
$ws.Cells.Item(297,11).Value = 500
$ws.Cells.Item(297,12).Value = 500
$ws.Cells.Item(297,13).Value = 500
$ws.Cells.Item(297,14).Value = '$/unidad'
$ws.Cells.Item(297,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(297,16).Value = 500
$ws.Cells.Item(297,17).Value = 1
$ws.Cells.Item(297,18).Value = 'Hortaliza'

# Row 298
$ws.Cells.Item(298,1).Value = 10
$ws.Cells.Item(298,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(298,3).Value = 'La Araucanía'
$ws.Cells.Item(298,4).Value = 44449
$ws.Cells.Item(298,5).Value = 9
$ws.Cells.Item(298,6).Value = 100112013
$ws.Cells.Item(298,7).Value = 'Alcachofa'
$ws.Cells.Item(298,8).Value = 'Española'
$ws.Cells.Item(298,9).Value = 'Primera'
$ws.Cells.Item(298,10).Value = 175
$ws.Cells.Item(298,11).Value = 14000
$ws.Cells.Item(298,12).Value = 15000
$ws.Cells.Item(298,13).Value = 14543
$ws.Cells.Item(298,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(298,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(298,16).Value = 485
$ws.Cells.Item(298,17).Value = 30
$ws.Cells.Item(298,18).Value = 'Hortaliza'

# Row 299
$ws.Cells.Item(299,1).Value = 10
$ws.Cells.Item(299,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(299,3).Value = 'La Araucanía'
$ws.Cells.Item(299,4).Value = 44449
$ws.Cells.Item(299,5).Value = 9
$ws.Cells.Item(299,6).Value = 100112013
$ws.Cells.Item(299,7).Value = 'Alcachofa'
$ws.Cells.Item(299,8).Value = 'Madrigal'
$ws.Cells.Item(299,9).Value = 'Primera'
$ws.Cells.Item(299,10).Value = 65
$ws.Cells.Item(299,11).Value = 14000
$ws.Cells.Item(299,12).Value = 14000
$ws.Cells.Item(299,13).Value = 14000
$ws.Cells.Item(299,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(299,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(299,16).Value = 350
$ws.Cells.Item(299,17).Value = 40
$ws.Cells.Item(299,18).Value = 'Hortaliza'

# Row 300
$ws.Cells.Item(300,1).Value = 10
$ws.Cells.Item(300,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(300,3).Value = 'La Araucanía'
$ws.Cells.Item(300,4).Value = 44837
$ws.Cells.Item(300,5).Value = 9
$ws.Cells.Item(300,6).Value = 100112013
$ws.Cells.Item(300,7).Value = 'Alcachofa'
$ws.Cells.Item(300,8).Value = 'Madrigal'
$ws.Cells.Item(300,9).Value = 'Primera'
$ws.Cells.Item(300,10).Value = 380
$ws.Cells.Item(300,11).Value = 11000
$ws.Cells.Item(300,12).Value = 12000
$ws.Cells.Item(300,13).Value = 11789
$ws.Cells.Item(300,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(300,15).Value = 'Región Metropolitana'
$ws.Cells.Item(300,16).Value = 295
$ws.Cells.Item(300,17).Value = 40
$ws.Cells.Item(300,18).Value = 'Hortaliza'

# Row 301
$ws.Cells.Item(301,1).Value = 10
$ws.Cells.Item(301,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(301,3).Value = 'La Araucanía'
$ws.Cells.Item(301,4).Value = 45049
$ws.Cells.Item(301,5).Value = 9
$ws.Cells.Item(301,6).Value = 100112013
$ws.Cells.Item(301,7).Value = 'Alcachofa'
$ws.Cells.Item(301,8).Value = 'Madrigal'
$ws.Cells.Item(301,9).Value = 'Primera'
$ws.Cells.Item(301,10).Value = 55
$ws.Cells.Item(301,11).Value = 18000
$ws.Cells.Item(301,12).Value = 18000
$ws.Cells.Item(301,13).Value = 18000
$ws.Cells.Item(301,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(301,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(301,16).Value = 450
$ws.Cells.Item(301,17).Value = 40
$ws.Cells.Item(301,18).Value = 'Hortaliza'

# Row 302
$ws.Cells.Item(302,1).Value = 10
$ws.Cells.Item(302,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(302,3).Value = 'La Araucanía'
$ws.Cells.Item(302,4).Value = 44462
$ws.Cells.Item(302,5).Value = 9
$ws.Cells.Item(302,6).Value = 100112013
$ws.Cells.Item(302,7).Value = 'Alcachofa'
$ws.Cells.Item(302,8).Value = 'Española'
$ws.Cells.Item(302,9).Value = 'Primera'
$ws.Cells.Item(302,10).Value = 50
$ws.Cells.Item(302,11).Value = 12000
$ws.Cells.Item(302,12).Value = 12000
$ws.Cells.Item(302,13).Value = 12000
$ws.Cells.Item(302,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(302,15).Value = 'Región Metropolitana'
$ws.Cells.Item(302,16).Value = 400
$ws.Cells.Item(302,17).Value = 30
$ws.Cells.Item(302,18).Value = 'Hortaliza'

# Row 303
$ws.Cells.Item(303,1).Value = 10
$ws.Cells.Item(303,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(303,3).Value = 'La Araucanía'
$ws.Cells.Item(303,4).Value = 44462
$ws.Cells.Item(303,5).Value = 9
$ws.Cells.Item(303,6).Value = 100112013
$ws.Cells.Item(303,7).Value = 'Alcachofa'
$ws.Cells.Item(303,8).Value = 'Madrigal'
$ws.Cells.Item(303,9).Value = 'Primera'
$ws.Cells.Item(303,10).Value = 150
$ws.Cells.Item(303,11).Value = 12000
$ws.Cells.Item(303,12).Value = 12000
$ws.Cells.Item(303,13).Value = 12000
$ws.Cells.Item(303,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(303,15).Value = 'Región Metropolitana'
$ws.Cells.Item(303,16).Value = 300
$ws.Cells.Item(303,17).Value = 40
$ws.Cells.Item(303,18).Value = 'Hortaliza'

# Row 304
$ws.Cells.Item(304,1).Value = 10
$ws.Cells.Item(304,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(304,3).Value = 'La Araucanía'
$ws.Cells.Item(304,4).Value = 44455
$ws.Cells.Item(304,5).Value = 9
$ws.Cells.Item(304,6).Value = 100112013
$ws.Cells.Item(304,7).Value = 'Alcachofa'
$ws.Cells.Item(304,8).Value = 'Argentina(o)'
$ws.Cells.Item(304,9).Value = 'Primera'
$ws.Cells.Item(304,10).Value = 80
$ws.Cells.Item(304,11).Value = 12000
$ws.Cells.Item(304,12).Value = 12000
$ws.Cells.Item(304,13).Value = 12000
$ws.Cells.Item(304,14).Value = '$/caja 40 unidades'
$ws.Cells.Item(304,15).Value = 'Región Metropolitana'
$ws.Cells.Item(304,16).Value = 300
$ws.Cells.Item(304,17).Value = 40
$ws.Cells.Item(304,18).Value = 'Hortaliza'

# Row 305
$ws.Cells.Item(305,1).Value = 10
$ws.Cells.Item(305,2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(305,3).Value = 'La Araucanía'
$ws.Cells.Item(305,4).Value = 44455
$ws.Cells.Item(305,5).Value = 9
$ws.Cells.Item(305,6).Value = 100112013
$ws.Cells.Item(305,7).Value = 'Alcachofa'
$ws.Cells.Item(305,8).Value = 'Española'
$ws.Cells.Item(305,9).Value = 'Primera'
$ws.Cells.Item(305,10).Value = 50
$ws.Cells.Item(305,11).Value = 12000
$ws.Cells.Item(305,12).Value = 12000
$ws.Cells.Item(305,13).Value = 12000
$ws.Cells.Item(305,14).Value = '$/caja 30 unidades'
$ws.Cells.Item(305,15).Value = 'Región Metropolitana'
$ws.Cells.Item(305,16).Value = 400
$ws.Cells.Item(305,17).Value = 30
$ws.Cells.Item(305,18).Value = 'Hortaliza'

# Ensure the date column (D) keeps the date number format for the newly appended rows
$dateFormat = $ws.Cells.Item(283,4).NumberFormat
$ws.Cells.Item(304,4).NumberFormat = $dateFormat
$ws.Cells.Item(305,4).NumberFormat = $dateFormat